$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 12000.5
$ws.Range("I62").Value = 2286.4285
$ws.Range("J62").Value = 34666.668
$ws.Range("K62").Value = 2286.4285
$ws.Range("L62").Value = 34666.668
$ws.Range("M62").Value = -1662.4285
$ws.Range("N62").Value = -35914.668
$ws.Range("H65").Value = 12000.5
$ws.Range("I65").Value = 2286.4285
$ws.Range("J65").Value = 34666.668
$ws.Range("K65").Value = 11432.1425
$ws.Range("L65").Value = 173333.34
$ws.Range("M65").Value = -8312.1425
$ws.Range("N65").Value = -179573.34
$ws.Range("H98").Value = 1891.6154
$ws.Range("I98").Value = 1235.909
$ws.Range("J98").Value = 5498
$ws.Range("K98").Value = 1235.909
$ws.Range("L98").Value = 5498
$ws.Range("M98").Value = 262.0909999999999
$ws.Range("N98").Value = -8494
$ws.Range("H122").Value = 1891.6154
$ws.Range("I122").Value = 1235.909
$ws.Range("J122").Value = 5498
$ws.Range("K122").Value = 3707.727
$ws.Range("L122").Value = 16494
$ws.Range("M122").Value = -1257.727
$ws.Range("N122").Value = -21394

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3400628.5
$ws.Range("I32").Value = 4867
$ws.Range("J32").Value = 20945396
$ws.Range("K32").Value = 4867
$ws.Range("L32").Value = 20945396
$ws.Range("M32").Value = -4580
$ws.Range("N32").Value = -20945970
$ws.Range("H45").Value = 4354
$ws.Range("I45").Value = 3512.3333
$ws.Range("J45").Value = 5300.875
$ws.Range("K45").Value = 3512.3333
$ws.Range("L45").Value = 5300.875
$ws.Range("M45").Value = -3135.3333
$ws.Range("N45").Value = -6054.875
$ws.Range("H61").Value = 1443.75
$ws.Range("I61").Value = 1443.75
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1443.75
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1231.75
$ws.Range("N61").Value = $null
$ws.Range("H122").Value = 2067.0952
$ws.Range("I122").Value = 1978.2142
$ws.Range("J122").Value = 2244.8572
$ws.Range("K122").Value = 5934.642599999999
$ws.Range("L122").Value = 6734.571599999999
$ws.Range("M122").Value = -3484.642599999999
$ws.Range("N122").Value = -11634.5716
$ws.Range("H136").Value = 1443.75
$ws.Range("I136").Value = 1443.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4331.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1781.25
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1422.425
$ws.Range("I58").Value = 1295.0646
$ws.Range("J58").Value = 1861.1111
$ws.Range("K58").Value = 1295.0646
$ws.Range("L58").Value = 1861.1111
$ws.Range("M58").Value = -1092.0646
$ws.Range("N58").Value = -2267.1111
$ws.Range("H62").Value = 2732.2632
$ws.Range("I62").Value = 2615.2144
$ws.Range("J62").Value = 3060
$ws.Range("K62").Value = 2615.2144
$ws.Range("L62").Value = 3060
$ws.Range("M62").Value = -1991.2144
$ws.Range("N62").Value = -4308
$ws.Range("H65").Value = 2732.2632
$ws.Range("I65").Value = 2615.2144
$ws.Range("J65").Value = 3060
$ws.Range("K65").Value = 13076.072
$ws.Range("L65").Value = 15300
$ws.Range("M65").Value = -9956.072
$ws.Range("N65").Value = -21540
$ws.Range("H122").Value = 3498.5715
$ws.Range("I122").Value = 1833.3334
$ws.Range("J122").Value = 4747.5
$ws.Range("K122").Value = 5500.0002
$ws.Range("L122").Value = 14242.5
$ws.Range("M122").Value = -3050.0002
$ws.Range("N122").Value = -19142.5
$ws.Range("H132").Value = 3075.5945
$ws.Range("I132").Value = 2606.1333
$ws.Range("J132").Value = 5087.5713
$ws.Range("K132").Value = 7818.3999
$ws.Range("L132").Value = 15262.7139
$ws.Range("M132").Value = -5288.3999
$ws.Range("N132").Value = -20322.7139
$ws.Range("H136").Value = 1422.425
$ws.Range("I136").Value = 1295.0646
$ws.Range("J136").Value = 1861.1111
$ws.Range("K136").Value = 3885.1938
$ws.Range("L136").Value = 5583.3333
$ws.Range("M136").Value = -1335.1938
$ws.Range("N136").Value = -10683.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 222514.55
$ws.Range("I5").Value = 248.83928
$ws.Range("J5").Value = 3334234.5
$ws.Range("K5").Value = 746.51784
$ws.Range("L5").Value = 10002703.5
$ws.Range("M5").Value = -634.51784
$ws.Range("N5").Value = -10002927.5
$ws.Range("H122").Value = 48864.434
$ws.Range("I122").Value = 385.4
$ws.Range("J122").Value = 55079.69
$ws.Range("K122").Value = 3468.6
$ws.Range("L122").Value = 495717.21
$ws.Range("M122").Value = -1018.6
$ws.Range("N122").Value = -500617.21
$ws.Range("H135").Value = 222514.55
$ws.Range("I135").Value = 248.83928
$ws.Range("J135").Value = 3334234.5
$ws.Range("K135").Value = 2239.55352
$ws.Range("L135").Value = 30008110.5
$ws.Range("M135").Value = 295.4464800000001
$ws.Range("N135").Value = -30013180.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2225.842
$ws.Range("I102").Value = 1693.1333
$ws.Range("J102").Value = 4223.5
$ws.Range("K102").Value = 1693.1333
$ws.Range("L102").Value = 4223.5
$ws.Range("M102").Value = -71.13329999999996
$ws.Range("N102").Value = -7467.5
$ws.Range("H122").Value = 2279.1667
$ws.Range("I122").Value = 1927.7778
$ws.Range("J122").Value = 3333.3333
$ws.Range("K122").Value = 5783.3334
$ws.Range("L122").Value = 9999.999899999999
$ws.Range("M122").Value = -3333.3334
$ws.Range("N122").Value = -14899.9999
$ws.Range("H126").Value = 10262.1875
$ws.Range("I126").Value = 2722.9412
$ws.Range("J126").Value = 18806.666
$ws.Range("K126").Value = 8168.823600000001
$ws.Range("L126").Value = 56419.99800000001
$ws.Range("M126").Value = -5698.823600000001
$ws.Range("N126").Value = -61359.99800000001
$ws.Range("H132").Value = 1645.1666
$ws.Range("I132").Value = 1323.3256
$ws.Range("J132").Value = 2459.2354
$ws.Range("K132").Value = 3969.976799999999
$ws.Range("L132").Value = 7377.706200000001
$ws.Range("M132").Value = -1439.976799999999
$ws.Range("N132").Value = -12437.7062

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3644.625
$ws.Range("I122").Value = 2336.375
$ws.Range("J122").Value = 4298.75
$ws.Range("K122").Value = 7009.125
$ws.Range("L122").Value = 12896.25
$ws.Range("M122").Value = -4559.125
$ws.Range("N122").Value = -17796.25
$ws.Range("H135").Value = 39199.6
$ws.Range("J135").Value = 39199.6
$ws.Range("L135").Value = 39199.6
$ws.Range("N135").Value = -49339.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 14286.471
$ws.Range("I122").Value = 17220.77
$ws.Range("J122").Value = 4750
$ws.Range("K122").Value = 51662.31
$ws.Range("L122").Value = 14250
$ws.Range("M122").Value = -49212.31
$ws.Range("N122").Value = -19150
